$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update curated dimension/measure labels in row 2
$ws.Range("A2").Value = "iaest-measure:grado"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:sexo"

# Update the "dim"/"medida" role labels in row 3
$ws.Range("A3").Value = "medida"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"

# Update the type labels in row 4
$ws.Range("A4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"

# Remove the mapping file row entirely (row 5: mapping-grado.xlsx / mapping-sexo.xlsx)
$ws.Rows.Item(5).Delete()
